# Adding Orders to Db implemented
# Update the "J" backlog column on Sheet1: mark the "Orders to Db" item as
# actioned (moves up), shift the remaining backlog entries up accordingly,
# drop the obsolete "Cost Of Rental..." entry, and append two freshly
# groomed backlog items at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J21").Value = "1) ORDERS - Add rental orders to the Db"
$ws.Range("J22").Value = "CarInfo - add for each individual car in list"
$ws.Range("J23").ClearContents()
$ws.Range("J24").Value = "Make checkout & confirmation pages dynamic"
$ws.Range("J25").Value = "Email With Confirmation - and PDF Invoice - confirmation GUID"
$ws.Range("J26").ClearContents()
$ws.Range("J27").ClearContents()
$ws.Range("J28").Value = "Add Address to Db - Edit from My Account"
$ws.Range("J29").Value = "Add Payment Details to DB - Edit from My Account"
$ws.Range("J30").ClearContents()
$ws.Range("J31").Value = "Checkout Discount code"
$ws.Range("J33").Value = "Notify user if no end date selected at checkout"

$ws.Range("J36").Select()
